$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 2 values, left to right (this order controls the order new
# entries land in xl/sharedStrings.xml).
$rowValues = @(
  @("A2", "1", $true),
  @("B2", "右线", $false),
  @("C2", "3", $true),
  @("D2", "192.168.155.3", $false),
  @("E2", "5820", $true),
  @("F2", "aasdasdf", $false),
  @("G2", "192.168.155.4", $false),
  @("H2", "aasdasdf", $false)
)

$scratch = $ws.Range("Z100")

foreach ($entry in $rowValues) {
    $addr = $entry[0]
    $text = $entry[1]
    $looksNumeric = $entry[2]

    if ($looksNumeric) {
        # Numeric-looking text ("1", "3", "5820") must stay text, not
        # become a number. Typing an apostrophe prefix forces text, but
        # it also stamps a quotePrefix cell style onto the cell. To land
        # the value as a plain, unstyled shared-string text cell, stage
        # it via a scratch cell, copy it, and paste-values-only into the
        # destination (paste-values drops the source formatting/style).
        $scratch.Value = "'" + $text
        $scratch.Copy()
        $ws.Range($addr).PasteSpecial(-4163)
    } else {
        $ws.Range($addr).Value = $text
    }
}

$scratch.Clear()
$excel.CutCopyMode = $false

